# 1/2 of the week02 task
# Update the triple-store style sheet: rename "peiper" subject to
# "peiper:Person", normalize predicates ("is"/"is " -> "is_a"), reorder
# the existing triples (officer now comes last among peiper's triples),
# and append two new triples (peiper born_in ..., himler knows peiper).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) stays the same: subject | predicate | object

# Row 2: peiper:Person | born | 30 January 1915
$ws.Range("A2").Value = "peiper:Person"
$ws.Range("B2").Value = "born"
$ws.Range("C2").Value = "30 January 1915"

# Row 3: peiper:Person | is_a | soldier
$ws.Range("A3").Value = "peiper:Person"
$ws.Range("B3").Value = "is_a"
$ws.Range("C3").Value = "soldier"

# Row 4: peiper:Person | die | 14 July 1976
$ws.Range("A4").Value = "peiper:Person"
$ws.Range("B4").Value = "die"
$ws.Range("C4").Value = "14 July 1976"

# Row 5: peiper:Person | is_a | nazi member
$ws.Range("A5").Value = "peiper:Person"
$ws.Range("B5").Value = "is_a"
$ws.Range("C5").Value = "nazi member"

# Row 6: peiper:Person | is_a | German
$ws.Range("A6").Value = "peiper:Person"
$ws.Range("B6").Value = "is_a"
$ws.Range("C6").Value = "German"

# Row 7: peiper:Person | is_a | officer  (new row, was previously row 2's data)
$ws.Range("A7").Value = "peiper:Person"
$ws.Range("B7").Value = "is_a"
$ws.Range("C7").Value = "officer"

# Row 8: peiper:Person | born_in | Wilmersdorf, Berlin, Germany  (new)
$ws.Range("A8").Value = "peiper:Person"
$ws.Range("B8").Value = "born_in"
$ws.Range("C8").Value = "Wilmersdorf, Berlin, Germany"

# Row 9: himler:Person | know | peiper  (new)
$ws.Range("A9").Value = "himler:Person"
$ws.Range("B9").Value = "know"
$ws.Range("C9").Value = "peiper"
